$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("C2").Value = 1000.0
$ws.Range("D2").Value = "15-01-2023"
$ws.Range("E2").Value = "Descripción del Activo 1"
$ws.Range("F2").Value = "Equipo de Oficina"
$ws.Range("G2").Value = "Nike"
$ws.Range("I2").Value = 66.67
$ws.Range("J2").Value = 933.33

# Row 3 updates
$ws.Range("A3").Value = 6.0
$ws.Range("B3").Value = "Monitor 2k"
$ws.Range("C3").Value = 1500.0
$ws.Range("D3").Value = "18-10-2023"
$ws.Range("E3").Value = "Monitor de alta resolucion"
$ws.Range("F3").Value = "Equipo de Oficina"
$ws.Range("G3").Value = "Puma"
$ws.Range("H3").Value = 10.0
$ws.Range("I3").Value = 0.0
$ws.Range("J3").Value = 1500.0
